#
# Insert a new "#elif TIZEN" directive line right before the existing
# "#elif WINDOWS" line in the preprocessor-directives table cell of the
# MAUI templates readme. This mirrors the commit's NuGet-template refresh,
# which added TIZEN to the platform list.
#
$d = $word.ActiveDocument

# Locate the paragraph that currently reads "#elif WINDOWS" (unique in doc).
$rng = $d.Content
$found = $rng.Find.Execute("#elif WINDOWS", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '#elif WINDOWS' paragraph to edit"
}

# Replace that paragraph with two paragraphs: the original run split into
# "#elif " + "TIZEN", followed by a brand-new paragraph (same formatting)
# that still reads "#elif WINDOWS".
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2D969A33" w14:textId="77777777" w:rsidR="00B33D13" w:rsidRPr="00B33D13" w:rsidRDefault="00B33D13" w:rsidP="00B33D13"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:pPr><w:r w:rsidRPr="00B33D13"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">#elif </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr><w:t>TIZEN</w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr><w:t>#elif WINDOWS</w:t></w:r></w:p>
'@

$rng.InsertXML($xml)
